# "updated in vivo data files"
#
# The ID->SampleType lookup sheet previously listed two full 24-sample
# cohorts (the "-2C" set in rows 2-25 and a duplicate "-1C" set in rows
# 26-49, plus a trailing blank row 50). The data was refreshed so that
# only the first cohort remains, and the old row 26 (KP1-1C / KP) was
# replaced with a single "inj" / "Stock" entry; everything below it
# (old rows 27-50) was removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old duplicate "-1C" cohort rows (27-49) and the trailing
# blank row (50), leaving the sheet with rows 1-26.
$ws.Range("A27:B50").EntireRow.Delete()

# Row 26 becomes a new "inj" / "Stock" record instead of "KP1-1C" / "KP".
$ws.Range("A26").Value = "inj"
$ws.Range("B26").Value = "Stock"

# Reflect the new selection used when the file was last saved.
$ws.Range("A26:XFD49").Select()
